# Daily attendance processing - 2026-02-02 06:05:02
#
# Normalize the ordering of the "Recorded By" (col G) audit-trail list for
# the rows whose list still has its legacy ordering: right-rotate the comma-
# separated list by one position (the last contributor moves to the front), e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"             -> "admin@admin.com, System"
#   "backup@backdoor.com, System, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (col G = "Recorded By") whose value needs the rotation applied,
# taken from the specific set of records touched by this processing run.
$rows = @(2, 3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 55, 58, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = $cell.Value2
    if ($current -eq $null) { continue }

    $parts = $current -split ', '
    $n = $parts.Count
    if ($n -lt 2) { continue }

    # Right-rotate by one: last entry moves to the front.
    $rotated = @($parts[$n - 1]) + $parts[0..($n - 2)]
    $cell.Value2 = [string]::Join(', ', $rotated)
}
